$wb = $excel.ActiveWorkbook

$wsMain    = $wb.Worksheets.Item("MAIN_CONTROLLER")
$wsData    = $wb.Worksheets.Item("DATASHEET")
$wsMobile  = $wb.Worksheets.Item("MOBILE_CONFIGURATION")

# ---------------------------------------------------------------------------
# 1) Rename the shared "FOS" label to "FOS_Till_UW" everywhere it now refers
#    to the new co-applicant ("Till UW") data sheet: MAIN_CONTROLLER!D2/E2,
#    DATASHEET!C3 and MOBILE_CONFIGURATION!C3.
# ---------------------------------------------------------------------------
$wsMain.Range("D2").Value = "FOS_Till_UW"
$wsMain.Range("E2").Value = "FOS_Till_UW"
$wsData.Range("C3").Value = "FOS_Till_UW"
$wsMobile.Range("C3").Value = "FOS_Till_UW"

# ---------------------------------------------------------------------------
# 2) MAIN_CONTROLLER A2:A6 become quote-prefixed text serial numbers
#    (1..5) instead of plain numbers.
# ---------------------------------------------------------------------------
$wsMain.Range("A2").Formula = "'1"
$wsMain.Range("A3").Formula = "'2"
$wsMain.Range("A4").Formula = "'3"
$wsMain.Range("A5").Formula = "'4"
$wsMain.Range("A6").Formula = "'5"

# ---------------------------------------------------------------------------
# 3) The old "FOSScroll.xlsx" file reference (DATASHEET!D3) becomes the new
#    co-applicant workbook name.
# ---------------------------------------------------------------------------
$wsData.Range("D3").Value = "FOS_Till_UW2.xlsx"

# ---------------------------------------------------------------------------
# 4) Drop the now-unused blank row that sat between row 6 and row 14 on
#    MAIN_CONTROLLER, shifting the trailing "N" marker up from H14 to H13.
# ---------------------------------------------------------------------------
$wsMain.Rows(7).Delete()

# ---------------------------------------------------------------------------
# 5) Update sheet selections / the active tab to match the saved view state.
#    Activate in order so MAIN_CONTROLLER ends up the active tab.
# ---------------------------------------------------------------------------
$wsMobile.Activate()
$wsMobile.Range("K3").Select()

$wsData.Activate()
$wsData.Range("D3").Select()

$wsMain.Activate()
$wsMain.Range("D2").Select()
